# Apply the updates described in the commit: refresh the "as of" date in the
# confidentiality disclosure text, and update the Weight / Percent Change
# figures for each holding row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; temporarily unprotect so the cells can be edited.
$wasProtected = $ws.ProtectContents
$ws.Unprotect()

# --- Update the disclosure date text -----------------------------------
$ws.Range("A38").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-13 for illustrative purposes only and are subject to change."

# --- Update Weight (D) and Percent Change (E) values --------------------
$updates = @(
    @{ Row = 2;  D = 0.03604454361702097;  E = -0.001575423395037467 },
    @{ Row = 3;  D = 0.02050345899335304;  E = 0.003142183817753219 },
    @{ Row = 4;  D = 0.0193892258601513;   E = 0.003642250101173472 },
    @{ Row = 5;  D = 0.03740713551615273;  E = 0.004305705059203468 },
    @{ Row = 6;  D = 0.03470265258364227;  E = 0.0004001600640257674 },
    @{ Row = 7;  D = 0.01992780973147699;  E = 0.003108606955507964 },
    @{ Row = 8;  D = 0.03707132720761505;  E = 0.006043600259011406 },
    @{ Row = 9;  D = 0.02045761767467795;  E = 0.005178052325581328 },
    @{ Row = 10; D = 0.02612645425840426;  E = 0.01758545741948225 },
    @{ Row = 11; D = 0.02375953482005643;  E = 0.0130363932645301 },
    @{ Row = 12; D = 0.05701490779816326;  E = 0.01390064063822072 },
    @{ Row = 13; D = 0.02496803174593916;  E = -0.002605135839225947 },
    @{ Row = 14; D = 0.02724419776285589;  E = 0.021286513362337 },
    @{ Row = 15; D = 0.03377142336339459;  E = 0.006848142441362892 },
    @{ Row = 16; D = 0.01987453468544918;  E = 0.01033005794910569 },
    @{ Row = 17; D = 0.03088868584373877;  E = 0.01947020974346114 },
    @{ Row = 18; D = 0.04204856903200558;  E = 0.001399580125962308 },
    @{ Row = 19; D = 0.126045764761608;    E = 0.002689979825151401 },
    @{ Row = 20; D = 0.00926118532692694;  E = 0.008846153846153726 },
    @{ Row = 21; D = 0.01524440662994526;  E = 0.01991181908690076 },
    @{ Row = 22; D = 0.01692349979574027;  E = 0.01897946484131907 },
    @{ Row = 23; D = 0.01547402620817364;  E = 0.00363636363636366 },
    @{ Row = 24; D = 0.02111023698851858;  E = 0.01219768664563636 },
    @{ Row = 25; D = 0.01237540085664932;  E = 0.01984348798211277 },
    @{ Row = 26; D = 0.04217014144808648;  E = 0.01542444774968321 },
    @{ Row = 27; D = 0.0242120629726531;   E = -0.0000980776775205694 },
    @{ Row = 28; D = 0.04547830498936552;  E = 0.006286266924564687 },
    @{ Row = 29; D = 0.05474919546033628;  E = 0.01646276104232114 },
    @{ Row = 30; D = 0.01269990370483807;  E = 0.01619433198380582 },
    @{ Row = 31; D = 0.02078010722170867;  E = 0.001155179052753308 },
    @{ Row = 32; D = 0.01328649706531673;  E = 0.004856726566294212 },
    @{ Row = 33; D = 0.04225377088080455;  E = 0.001033591731266048 },
    @{ Row = 34; D = 0.01673538519523123;  E = 0.008162031438936124 },
    @{ Row = 35; D = $null;                E = 0.007733408778782813 }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# Restore sheet protection to its original (unprotected-from-the-password
# perspective, but protected) state.
if ($wasProtected) {
    $ws.Protect()
}
